# Update "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
# (shared with de-de "Correspond Handoff Datetime" for the same file).
$wsOverview.Range("G2").Value = "2016-08-18 01:01:28"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file.
$wsZhCn.Range("H2").Value = "2016-08-18 01:01:23"
$wsZhCn.Range("K2").Value = "2016-08-18 01:01:49"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview G2) and
# "Correspond Handback DateTime" for the first file.
$wsDeDe.Range("H2").Value = "2016-08-18 01:01:28"
$wsDeDe.Range("K2").Value = "2016-08-18 01:01:56"
